$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update all order statuses (column H, rows 3-6) to "En attente"
$ws.Range("H3").Value = "En attente"
$ws.Range("H4").Value = "En attente"
$ws.Range("H5").Value = "En attente"
$ws.Range("H6").Value = "En attente"

# Move the active selection, matching the recorded cursor position after the edit
$ws.Range("E12").Select()
